$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45171 -> 45172) for every data row (rows 2 through 480).
$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
if ($lastRow -lt 480) { $lastRow = 480 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
